$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 5, shifting existing rows 5-17 down to 6-18
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new exposure site entry
$ws.Range("A5").Value = "Lakes Entrance"
$ws.Range("B5").Value = "Chants Summer Carnival - Footbridge, Lakes Entrance VIC 3909"
$ws.Range("C5").Value = "29/12/2020 7:00pm-9:30pm"
$ws.Range("D5").Value = "Case attended carnival"

# Fix the date for the European Bier Cafe City row (now row 13, was row 12)
$ws.Range("C13").Value = "28/12/2020 8:00pm-9:30pm"

# Fix typo in the Rockpool Bar and Grill exposure date (now row 18, was row 17)
$ws.Range("C18").Value = "23/12/2020 1:00pm-1:30pm"
